$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting existing rows 207-221 down to 208-222
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new data record
$ws.Range("A207").Value = 7
$ws.Range("B207").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C207").Value = "Ñuble"
$ws.Range("D207").Value = 45021
$ws.Range("E207").Value = 16
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100101
$ws.Range("H207").Value = "Berries"
$ws.Range("I207").Value = 100101007
$ws.Range("J207").Value = "Kiwi"
$ws.Range("K207").Value = "Hayward"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 60
$ws.Range("N207").Value = 16000
$ws.Range("O207").Value = 16000
$ws.Range("P207").Value = 16000
$ws.Range("Q207").Value = "$/bandeja 18 kilos"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 889
$ws.Range("T207").Value = 18
